# addBookCommand, update some test partially
#
# Renames the class-diagram boxes from the AddressBook/Person sample
# domain to the Book/BookShelf domain, and nudges the resized
# "Address" -> "ReviewList" box (plus its connector) to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# --- Simple label renames (class / attribute boxes) ---------------------

(Get-ShapeById $s 46).TextFrame.TextRange.Text = "VersionedBookShelf"
(Get-ShapeById $s 49).TextFrame.TextRange.Text = "UniqueBookList"
(Get-ShapeById $s 62).TextFrame.TextRange.Text = "Book"
(Get-ShapeById $s 80).TextFrame.TextRange.Text = "Author"
(Get-ShapeById $s 83).TextFrame.TextRange.Text = "Rating"
(Get-ShapeById $s 55).TextFrame.TextRange.Text = "BookShelf"

# "<<interface>>ReadOnlyAddressBook" -> "<<interface>>ReadOnlyBookShelf"
# (only the second run, after the line break, changes)
$readOnly = Get-ShapeById $s 100
$readOnly.TextFrame.TextRange.Runs(2).Text = "ReadOnlyBookShelf"

# --- "Address" -> "ReviewList" box, which also grew wider ---------------

$address = Get-ShapeById $s 85
$address.TextFrame.TextRange.Text = "ReviewList"
$address.Left = 607.2753
$address.Width = 63.97787401574803

# The elbow connector feeding into that box shrinks slightly to follow
# the box's new position.
$addressConnector = Get-ShapeById $s 86
$addressConnector.Width = 34.2048031
